$d = $word.ActiveDocument

function Set-ParagraphBold($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # Exclude trailing paragraph mark from the range
    $r2 = $d.Range($r.Start, $r.End - 1)
    if ($newText) {
        $r2.Text = $newText
    }
    $r2.Font.Bold = 1
}

# 1. "Nestor Wilke" -> bold (no text change)
Set-ParagraphBold 1 $null

# 2. "Expérience professionnelle" -> bold (no text change)
Set-ParagraphBold 5 $null

# 3. "Responsable de l’équipe d’animation" -> bold + text change
Set-ParagraphBold 6 "Chef d’équipe d’animation"

# 4. "Concepteur d’animation principal" (2nd job heading) -> bold + text change
Set-ParagraphBold 12 "Concepteur d’animation senior"

# 5. "Concepteur d’animation" (3rd job heading) -> bold (no text change)
Set-ParagraphBold 18 $null

# 6. "Licence en Beaux-Arts spécialisée en animation" -> bold + text change
Set-ParagraphBold 24 "Licence d’arts plastiques en animation"
